$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing table (Table2) by one row so the table ref/autofilter
# and sheet dimension grow from H24 to H25, matching the source file.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Populate the new row (row 25) with the task data. The order in which new
# shared strings are first introduced matters (they are appended to
# sharedStrings.xml in first-seen order), so the "End date" text is written
# before the "Task" text to reproduce the original author's string order.
$ws.Cells.Item(25, 1).Value = 24

# The "End date" column for this row stores the literal text "2025-08-24"
# (not a real date value) while still using the date-formatted style, just
# like the rows above it (F23, F24). Entering the text directly causes the
# host to auto-parse it into a date serial, so we stage the literal text (via
# a leading apostrophe to force text) in a scratch cell, copy only the value
# across, then copy the number format from the row above so the cell's style
# matches the existing date-styled column exactly.
$scratch = $ws.Cells.Item(1048576, 16384)
$scratch.Value = "'2025-08-24"
$scratch.Copy() | Out-Null
$ws.Cells.Item(25, 6).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(24, 6).Copy() | Out-Null
$ws.Cells.Item(25, 6).PasteSpecial(-4122) | Out-Null
$scratch.Clear() | Out-Null

$ws.Cells.Item(25, 2).Value = "Submit the SRS"
$ws.Cells.Item(25, 3).Value = "Not Started"
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 45860
$ws.Cells.Item(25, 7).Value = "Aishwarrya VP"

$excel.CutCopyMode = 0

# Match the author's final selection.
$ws.Range("E29").Select() | Out-Null
